$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Product Name" column (column B) -----------------------
# This shifts Date / Total Amount / Payment Method left from C/D/E into
# B/C/D, giving the new 4-column layout (Order ID, Date, Total Amount,
# Payment Method).
$ws.Columns("B").Delete()

# --- Clear all remaining old data so we can lay out the fresh report ---
$ws.Range("A1:D10").ClearContents()

# Helper: write a value that must stay literal text even though it looks
# like a number (e.g. "87945.00") so it keeps its exact formatting and is
# stored as a shared string instead of being coerced into a numeric cell.
# Building it via a text formula in a scratch cell and pasting back as a
# value keeps the destination cell on the workbook's default style (no
# NumberFormat side effects).
function Set-TextValue($targetAddr, $value) {
    $helper = $ws.Range("Z1")
    $helper.Formula = '="' + $value + '"'
    $helper.Copy()
    $ws.Range($targetAddr).PasteSpecial(-4163)
    $helper.ClearContents()
}

# --- Header row ----------------------------------------------------------
$ws.Range("A1").Value = "Order ID"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Total Amount"
$ws.Range("D1").Value = "Payment Method"

# --- Order rows ------------------------------------------------------------
$orders = @(
    @("DPRJV1CL", "23/1/2024", "87945.00", "COD"),
    @("WVAJ27DU", "24/1/2024", "108945.00", "COD"),
    @("4A690OVT", "24/1/2024", "87945.00", "COD"),
    @("BA8WTN90", "24/1/2024", "24500.00", "COD"),
    @("05ZJMUH8", "24/1/2024", "21000.00", "COD"),
    @("136J6JJT", "25/1/2024", "320442.00", "COD")
)

$row = 2
foreach ($order in $orders) {
    $ws.Range("A$row").Value = $order[0]
    $ws.Range("B$row").Value = $order[1]
    Set-TextValue "C$row" $order[2]
    $ws.Range("D$row").Value = $order[3]
    $row++
}

# --- Totals row (now row 8, in columns C:D) -------------------------------
$ws.Range("C8").Value = "Total Sales Amount"
Set-TextValue "D8" "650777.00"
